$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" (sheet1): row 3 is the d6c199f8 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 08:57:36"

# --- Sheet "zh-cn" (sheet2): row 3 is the d6c199f8 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-13 08:57:29"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/da4d932e5e8bc512d099078b2bd4d749e71577a4/e2e/d6c199f8-1caa-40ef-950b-3ab0e6e69777.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/354de7c5164b06aaec8aebd6ead0dabac663d83f/e2e/d6c199f8-1caa-40ef-950b-3ab0e6e69777.md."
$wsZhCn.Range("P:P").ColumnWidth = 39.1666666666667

# --- Sheet "de-de" (sheet3): row 3 is the d6c199f8 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-13 08:57:36"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/da4d932e5e8bc512d099078b2bd4d749e71577a4/e2e/d6c199f8-1caa-40ef-950b-3ab0e6e69777.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/354de7c5164b06aaec8aebd6ead0dabac663d83f/e2e/d6c199f8-1caa-40ef-950b-3ab0e6e69777.md."
$wsDeDe.Range("P:P").ColumnWidth = 39.1666666666667
